$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Create the new "TotalToTarget (lbs)" column header in G1 first (so its shared
#    string is appended before the E1/F1 renames, matching the target string order).
#    Copy F1's number format first so G1 picks up the same header style (s="2").
$ws.Range("F1").Copy()
$ws.Range("G1").PasteSpecial(-4122)
$ws.Range("G1").Value = "TotalToTarget (lbs)"

# 2. Rename the existing headers.
$ws.Range("E1").Value = "Weight (lbs)"
$ws.Range("F1").Value = "Gain/Loss Amount (lbs)"

# 3. Add the TotalToTarget formula for the first data row.
$ws.Range("G2").Formula = "=E2-210"

# Widen F to fit its longer header text and size the new G column to fit its header.
$ws.Columns("F").ColumnWidth = 19.666
$ws.Columns("G").ColumnWidth = 15.666

# 4. Bump the second measurement's date forward a week.
$ws.Range("B3").Value = 43101

# 5. Drop the trailing placeholder rows that only held WeightID/Date.
$ws.Rows("4:8").Delete()

# 6. Leave the selection where the user left it after entering the new formula.
$ws.Range("E2").Select() | Out-Null
